$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,14
$data[0,0] = 13.1020553966038
$data[0,1] = 11.97891613822523
$data[0,2] = 0
$data[0,3] = 16.86610535988619
$data[0,4] = 37.57566155749477
$data[0,5] = 31.21283244510076
$data[0,6] = 15.03816081369865
$data[0,7] = 0
$data[0,8] = 7.874477370322277
$data[0,9] = 8.589507243371235
$data[0,10] = 12.45215293674379
$data[0,11] = 0
$data[0,12] = 19.28880467906017
$data[0,13] = 23.19134679964884
$data[1,0] = 12.86603731740907
$data[1,1] = 12.0126186248821
$data[1,2] = 0
$data[1,3] = 16.86445353970958
$data[1,4] = 37.59597018204517
$data[1,5] = 31.31377107261725
$data[1,6] = 15.08333806467516
$data[1,7] = 0
$data[1,8] = 7.864912761074475
$data[1,9] = 8.418696921009069
$data[1,10] = 12.42957505037062
$data[1,11] = 0
$data[1,12] = 19.34239775978782
$data[1,13] = 23.2690996495586
$data[2,0] = 12.72097775379347
$data[2,1] = 12.03438854475512
$data[2,2] = 0
$data[2,3] = 16.8661099469756
$data[2,4] = 37.6165401934628
$data[2,5] = 31.38369401309318
$data[2,6] = 15.11304127174312
$data[2,7] = 0
$data[2,8] = 7.858973889705817
$data[2,9] = 8.313054386460829
$data[2,10] = 12.417455128145
$data[2,11] = 0
$data[2,12] = 19.37691414553237
$data[2,13] = 23.32085908786378
$data[3,0] = 12.66190613727701
$data[3,1] = 12.04353140158601
$data[3,2] = 0
$data[3,3] = 16.8674583619831
$data[3,4] = 37.62696028808472
$data[3,5] = 31.41418054170913
$data[3,6] = 15.12563991759136
$data[3,7] = 0
$data[3,8] = 7.856537323947144
$data[3,9] = 8.269868043485086
$data[3,10] = 12.41295814876736
$data[3,11] = 0
$data[3,12] = 19.39138580378537
$data[3,13] = 23.3429615232065
$data[4,0] = 12.65210212336976
$data[4,1] = 12.04506598237294
$data[4,2] = 0
$data[4,3] = 16.86772298720618
$data[4,4] = 37.62881360947173
$data[4,5] = 31.4193629654603
$data[4,6] = 15.12776178346603
$data[4,7] = 0
$data[4,8] = 7.856131747930248
$data[4,9] = 8.262690429207126
$data[4,10] = 12.41223822615766
$data[4,11] = 0
$data[4,12] = 19.39381336243606
$data[4,13] = 23.34669260416739
$data[5,0] = 12.72018082275552
$data[5,1] = 12.03451074847277
$data[5,2] = 0
$data[5,3] = 16.86612540326587
$data[5,4] = 37.61667247184617
$data[5,5] = 31.38409710594647
$data[5,6] = 15.11320917930933
$data[5,7] = 0
$data[5,8] = 7.858941095568658
$data[5,9] = 8.312472436415595
$data[5,10] = 12.41739268604309
$data[5,11] = 0
$data[5,12] = 19.37710767023619
$data[5,13] = 23.32115308012051
$data[6,0] = 13.02075217010212
$data[6,1] = 11.99031379197572
$data[6,2] = 0
$data[6,3] = 16.86498271507065
$data[6,4] = 37.58098342087379
$data[6,5] = 31.24598346491044
$data[6,6] = 15.05333052377948
$data[6,7] = 0
$data[6,8] = 7.871193350128373
$data[6,9] = 8.53080358147357
$data[6,10] = 12.44400862405434
$data[6,11] = 0
$data[6,12] = 19.30695007772381
$data[6,13] = 23.21732154893299
$data[7,0] = 13.60571129122859
$data[7,1] = 11.91215027848029
$data[7,2] = 0
$data[7,3] = 16.88384101783797
$data[7,4] = 37.57520439235497
$data[7,5] = 31.03843850543999
$data[7,6] = 14.95147402660877
$data[7,7] = 0
$data[7,8] = 7.894686454367025
$data[7,9] = 8.950514606132087
$data[7,10] = 12.50985886622655
$data[7,11] = 0
$data[7,12] = 19.18209522421619
$data[7,13] = 23.04561966945879
$data[8,0] = 14.02846349737075
$data[8,1] = 11.85986062376526
$data[8,2] = 0
$data[8,3] = 16.91042278403707
$data[8,4] = 37.60996086684302
$data[8,5] = 30.92485876887776
$data[8,6] = 14.88610203180116
$data[8,7] = 0
$data[8,8] = 7.911611057529091
$data[8,9] = 9.250709035670729
$data[8,10] = 12.56631629077256
$data[8,11] = 0
$data[8,12] = 19.09804921144042
$data[8,13] = 22.93895579108613
$data[9,0] = 14.21840430683613
$data[9,1] = 11.83717780861062
$data[9,2] = 0
$data[9,3] = 16.9252444004071
$data[9,4] = 37.63419185010999
$data[9,5] = 30.88169628948508
$data[9,6] = 14.85841160752738
$data[9,7] = 0
$data[9,8] = 7.919235312133297
$data[9,9] = 9.384910755366633
$data[9,10] = 12.59369562545737
$data[9,11] = 0
$data[9,12] = 19.06146789716806
$data[9,13] = 22.89466931959242
$data[10,0] = 14.28992188930054
$data[10,1] = 11.82874639930071
$data[10,2] = 0
$data[10,3] = 16.93124601754108
$data[10,4] = 37.64457234323851
$data[10,5] = 30.86657915888732
$data[10,6] = 14.84822002338145
$data[10,7] = 0
$data[10,8] = 7.922111427024327
$data[10,9] = 9.435344745461778
$data[10,10] = 12.60430201924054
$data[10,11] = 0
$data[10,12] = 19.0478518865155
$data[10,13] = 22.87850882413962
$data[11,0] = 14.27453852606482
$data[11,1] = 11.83055523415614
$data[11,2] = 0
$data[11,3] = 16.92993622005693
$data[11,4] = 37.64228324412127
$data[11,5] = 30.86978024063549
$data[11,6] = 14.85040188387469
$data[11,7] = 0
$data[11,8] = 7.921492498813048
$data[11,9] = 9.424500723489434
$data[11,10] = 12.60200722893517
$data[11,11] = 0
$data[11,12] = 19.05077383311512
$data[11,13] = 22.88196214148125
$data[12,0] = 14.22429667204786
$data[12,1] = 11.83648098716469
$data[12,2] = 0
$data[12,3] = 16.92573038210654
$data[12,4] = 37.6350217529119
$data[12,5] = 30.88042796777149
$data[12,6] = 14.85756724549053
$data[12,7] = 0
$data[12,8] = 7.919472154321309
$data[12,9] = 9.389067966019324
$data[12,10] = 12.59456347387411
$data[12,11] = 0
$data[12,12] = 19.06034296466055
$data[12,13] = 22.89332755962658
$data[13,0] = 14.1934668838276
$data[13,1] = 11.84013124881861
$data[13,2] = 0
$data[13,3] = 16.9232047346193
$data[13,4] = 37.63073057633108
$data[13,5] = 30.88710999848353
$data[13,6] = 14.86199453931926
$data[13,7] = 0
$data[13,8] = 7.918233191299727
$data[13,9] = 9.367312882061048
$data[13,10] = 12.59003484158472
$data[13,11] = 0
$data[13,12] = 19.0662351100019
$data[13,13] = 22.90036864699037
$data[14,0] = 14.01599681753669
$data[14,1] = 11.86136515675519
$data[14,2] = 0
$data[14,3] = 16.90950878736424
$data[14,4] = 37.60854625462562
$data[14,5] = 30.92785108427379
$data[14,6] = 14.88795284654005
$data[14,7] = 0
$data[14,8] = 7.91111126628358
$data[14,9] = 9.241887170572683
$data[14,10] = 12.5645606592986
$data[14,11] = 0
$data[14,12] = 19.10047303881956
$data[14,13] = 22.94193531123967
$data[15,0] = 13.90646790812417
$data[15,1] = 11.87467375282711
$data[15,2] = 0
$data[15,3] = 16.90180344399165
$data[15,4] = 37.59708949492999
$data[15,5] = 30.95502621516085
$data[15,6] = 14.90440168002091
$data[15,7] = 0
$data[15,8] = 7.90672284332366
$data[15,9] = 9.164304788162982
$data[15,10] = 12.54936375312732
$data[15,11] = 0
$data[15,12] = 19.12189923324851
$data[15,13] = 22.96852047493083
$data[16,0] = 13.84325042069785
$data[16,1] = 11.88243245995933
$data[16,2] = 0
$data[16,3] = 16.89762866924857
$data[16,4] = 37.591292835138
$data[16,5] = 30.97145697023096
$data[16,6] = 14.91405536147981
$data[16,7] = 0
$data[16,8] = 7.904191711050674
$data[16,9] = 9.119462066167737
$data[16,10] = 12.5407829314641
$data[16,11] = 0
$data[16,12] = 19.13437853047352
$data[16,13] = 22.9842101743866
$data[17,0] = 13.82181051793656
$data[17,1] = 11.88507729987271
$data[17,2] = 0
$data[17,3] = 16.89625942732267
$data[17,4] = 37.58946654607057
$data[17,5] = 30.977157456723
$data[17,6] = 14.91735704500529
$data[17,7] = 0
$data[17,8] = 7.903333514145981
$data[17,9] = 9.104242865339019
$data[17,10] = 12.53790526343276
$data[17,11] = 0
$data[17,12] = 19.13863054586296
$data[17,13] = 22.9895908758499
$data[18,0] = 13.91815065310693
$data[18,1] = 11.87324627604948
$data[18,2] = 0
$data[18,3] = 16.90259710287861
$data[18,4] = 37.59822705823544
$data[18,5] = 30.95205051017253
$data[18,6] = 14.90263072693361
$data[18,7] = 0
$data[18,8] = 7.907190727651743
$data[18,9] = 9.172586614690967
$data[18,10] = 12.5509649657678
$data[18,11] = 0
$data[18,12] = 19.11960228988113
$data[18,13] = 22.9656491778772
$data[19,0] = 14.23906554257558
$data[19,1] = 11.83473616494982
$data[19,2] = 0
$data[19,3] = 16.92695521079035
$data[19,4] = 37.63712198572829
$data[19,5] = 30.87726712298825
$data[19,6] = 14.85545462326
$data[19,7] = 0
$data[19,8] = 7.920065879812934
$data[19,9] = 9.399486236025654
$data[19,10] = 12.59674345949117
$data[19,11] = 0
$data[19,12] = 19.05752586678686
$data[19,13] = 22.88997270411489
$data[20,0] = 14.44638978576598
$data[20,1] = 11.81048868694463
$data[20,2] = 0
$data[20,3] = 16.94514016777044
$data[20,4] = 37.66956072428206
$data[20,5] = 30.83554879553051
$data[20,6] = 14.82633689203727
$data[20,7] = 0
$data[20,8] = 7.928416148731217
$data[20,9] = 9.545513136534296
$data[20,10] = 12.62804938353793
$data[20,11] = 0
$data[20,12] = 19.01833362818649
$data[20,13] = 22.84406865627592
$data[21,0] = 14.33597914426016
$data[21,1] = 11.8233459622882
$data[21,2] = 0
$data[21,3] = 16.93522844669732
$data[21,4] = 37.65160761373691
$data[21,5] = 30.85715842763381
$data[21,6] = 14.84172077518021
$data[21,7] = 0
$data[21,8] = 7.92396543301874
$data[21,9] = 9.467797593326791
$data[21,10] = 12.61121580182077
$data[21,11] = 0
$data[21,12] = 19.0391254651228
$data[21,13] = 22.86824298556821
$data[22,0] = 13.91286965091431
$data[22,1] = 11.87389130365499
$data[22,2] = 0
$data[22,3] = 16.90223749473869
$data[22,4] = 37.59771030469665
$data[22,5] = 30.95339331172126
$data[22,6] = 14.90343076107269
$data[22,7] = 0
$data[22,8] = 7.90697922257422
$data[22,9] = 9.168843143600307
$data[22,10] = 12.55024057094127
$data[22,11] = 0
$data[22,12] = 19.1206402361068
$data[22,13] = 22.96694602811921
$data[23,0] = 13.44839562435748
$data[23,1] = 11.93238995731016
$data[23,2] = 0
$data[23,3] = 16.87649346842417
$data[23,4] = 37.56990455721503
$data[23,5] = 31.08777451724202
$data[23,6] = 14.97736536335896
$data[23,7] = 0
$data[23,8] = 7.888389131748314
$data[23,9] = 8.83819634894896
$data[23,10] = 12.49060709620025
$data[23,11] = 0
$data[23,12] = 19.21451707264341
$data[23,13] = 23.08864994572665

$ws.Range("B2:O25").Value2 = $data
